$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row (19) of data ---
$ws.Range("A19").Value = "Number of 1 Bits"
$ws.Range("B19").Value = "Binary"
$ws.Range("C19").Value = "Yes"
$ws.Range("D19").Value = "No"
$ws.Range("E19").Value = "Easy"
$ws.Range("F19").Value = "Easy"
$ws.Range("G19").Value = "191 - Number of 1 Bits"

# --- Add the hyperlink for the new problem link cell ---
$ws.Hyperlinks.Add($ws.Cells.Item(19, 7), "191 - Number of 1 Bits")

# Re-apply the same "Hyperlink" cell style used by the other rows in column G
# (Hyperlinks.Add() stamps its own fresh style; put it back in line with G2:G18)
$ws.Range("G18").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Extend the data validation ranges to include row 19 ---
# The COM model doesn't expose an in-place "resize" for Validation, so each
# rule is deleted and re-added against the new, larger range (this keeps the
# rule count / ordering the same as before, just with the bigger sqref).
$ws.Range("E2:F18").Validation.Delete()
$ws.Range("C2:C18").Validation.Delete()
$ws.Range("B2:B18").Validation.Delete()
$ws.Range("D2:D18").Validation.Delete()

$ws.Range("E2:F19").Validation.Add(3, 1, 1, '"Easy, Medium, Hard"')
$ws.Range("E2:F19").Validation.IgnoreBlank = $true
$ws.Range("E2:F19").Validation.InCellDropdown = $true

$ws.Range("C2:C19").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("C2:C19").Validation.InCellDropdown = $true

$ws.Range("B2:B19").Validation.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap"')
$ws.Range("B2:B19").Validation.IgnoreBlank = $true
$ws.Range("B2:B19").Validation.InCellDropdown = $true

$ws.Range("D2:D19").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("D2:D19").Validation.IgnoreBlank = $true
$ws.Range("D2:D19").Validation.InCellDropdown = $true

# --- Update the selected cell shown in the sheet view ---
$ws.Range("F30").Select()
